# "added get ranges from sheet functionality"
#
# Replace the single workbook-scoped "InstrumentIdentifier" defined name
# (Instruments!$A$2:$A$6) with sheet-scoped defined names that let each
# instrument sheet resolve its own identifier/class cell, plus a
# workbook-level "Instruments" name that keeps the original range.

$wb = $excel.ActiveWorkbook

$wsInstruments = $wb.Worksheets.Item("Instruments")
$wsCAN1        = $wb.Worksheets.Item("CAN1")
$wsPCAN        = $wb.Worksheets.Item("PCAN_USBBUS1  0x51")

# Drop the old single workbook-scoped name.
$wb.Names.Item("InstrumentIdentifier").Delete()

# Per-sheet "InstrumentClass" / "InstrumentIdentifier" names.
$wsCAN1.Names.Add("InstrumentClass", "='CAN1'!`$B`$2")
$wsPCAN.Names.Add("InstrumentClass", "='PCAN_USBBUS1  0x51'!`$B`$3")
$wsCAN1.Names.Add("InstrumentIdentifier", "='CAN1'!`$A`$2")
$wsPCAN.Names.Add("InstrumentIdentifier", "='PCAN_USBBUS1  0x51'!`$A`$3")

# Workbook range of instruments, now scoped to the Instruments sheet.
$wsInstruments.Names.Add("Instruments", "=Instruments!`$A`$2:`$A`$6")

# Rename the header label on the Instruments sheet.
$wsInstruments.Range("A1").Value = "Study Instruments"

# Update remembered selections: CAN1!A2, then Instruments!A8 last so
# Instruments keeps being the active/tabSelected sheet.
$wsCAN1.Activate()
$wsCAN1.Range("A2").Select()

$wsInstruments.Activate()
$wsInstruments.Range("A8").Select()
